$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 (Ano = 2025) with refreshed data values
$ws.Range("B7").Value = 2765916.02
$ws.Range("C7").Value = -37.74775177709457
$ws.Range("D7").Value = 2816
$ws.Range("E7").Value = 2816
$ws.Range("F7").Value = 982.2144957386364
$ws.Range("G7").Value = 4.696962920340941
